# Adds the new "2022-12-22" (serial 44942) snapshot row to the "Data" table,
# mirroring the upstream check_ssl_cert_stats.xlsx commit that appended one more
# day of stats and removed the (now OS-OpenSSL-version-dependent) --security-level
# test results from the still-running GH workflow columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row (this shifts ref/autoFilter to A1:AJ59 and
# dimension to match).
$newListRow = $lo.ListRows.Add()

# Carry the previous row's number formats down onto the freshly inserted
# row before we populate it, so date/int/delta columns keep their styling.
$ws.Range("A58:AJ58").Copy()
$ws.Range("A59:AJ59").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 59 values -------------------------------------------------------
$ws.Range("A59").Value2 = 44942
$ws.Range("B59").Value2 = 326
$ws.Range("C59").Value2 = 122
$ws.Range("D59").Value2 = 110
$ws.Range("E59").Value2 = 263
$ws.Range("F59").Value2 = 216
$ws.Range("G59").Value2 = 5491
$ws.Range("H59").Formula = "=Data[[#This Row],[LoC]]-G58"
$ws.Range("I59").Value2 = 6364
$ws.Range("J59").Value2 = 1916
$ws.Range("K59").Value2 = 265
$ws.Range("L59").Value2 = 285
$ws.Range("M59").Value2 = 106
$ws.Range("N59").Value2 = 60
$ws.Range("O59").Formula = "=SUM(Data[[#This Row],[Shell]:[Bash]])"
$ws.Range("P59").Formula = "=Data[[#This Row],[Total]]-O58"
$ws.Range("Q59").Value2 = 2002
$ws.Range("R59").Value2 = 4201
$ws.Range("S59").Value2 = 66282
$ws.Range("T59").Value2 = 46002
$ws.Range("U59").Value2 = 0
$ws.Range("V59").Value2 = 0
$ws.Range("W59").Value2 = 247
$ws.Range("X59").Formula = "=Data[[#This Row],[Open issues]]+Data[[#This Row],[Closed issues]]"
$ws.Range("Y59").Value2 = 0
$ws.Range("Z59").Value2 = 170
$ws.Range("AA59").Formula = "=Data[[#This Row],[Open pull requests]]+Data[[#This Row],[Closed pull requests]]"
$ws.Range("AB59").Value2 = 156
# AC59:AI59 (Tests / GH workflows / Running / Failed / OK / Cancelled / GH runs)
# intentionally stay blank: the author dropped the in-flight --security-level
# run's numbers for this row.
$ws.Range("AC59").ClearContents()
$ws.Range("AD59").ClearContents()
$ws.Range("AE59").ClearContents()
$ws.Range("AF59").ClearContents()
$ws.Range("AG59").ClearContents()
$ws.Range("AH59").ClearContents()
$ws.Range("AI59").ClearContents()
$ws.Range("AJ59").Formula = "=SUM(Data[[#This Row],[Running]:[GH runs]])"

# The Stars/Forks "latest value" highlight style moves from the old last
# row down to the new one.
$ws.Range("B58:C58").Copy()
$ws.Range("B59:C59").PasteSpecial(-4122)
$ws.Range("B58:C58").ClearFormats()
$ws.Range("B58").Value2 = 325
$ws.Range("C58").Value2 = 121

# --- View state ------------------------------------------------------------
# Freeze pane follows the growing table (was row 48, now row 52) and the
# active selection follows the last populated data cell.
$ws.Range("A52").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("AI59").Select()
$ws.Activate()
